{"js": "// The document contains one table of structural-variant calls. Two kinds\n// of cleanup are applied to its data cells:\n//\n//   1. Numeric values that were serialized with a trailing \".000\" (e.g.\n//      \"137174039.000\") are rewritten without that suffix (\"137174039\").\n//   2. A few \"Event Size\" cells that were left blank (no <w:t> run at all)\n//      get the literal text \"<NA>\" inserted.\n\n// --- 1. Strip the \".000\" suffix from every matching numeric run -----------\n// Wildcard pattern: one-or-more digits immediately followed by \".000\".\nconst body = context.document.body;\nconst trailingZeros = body.search(\"[0-9]@.000\", { matchWildcards: true });\ntrailingZeros.load(\"items/text\");\nawait context.sync();\n\nfor (const found of trailingZeros.items) {\n  const trimmed = found.text.slice(0, -4); // drop the literal \".000\"\n  found.insertText(trimmed, Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// --- 2. Fill blank \"Event Size\" cells with the text \"<NA>\" -----------------\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount,values\");\nawait context.sync();\n\n// Locate the \"Event Size\" column by its header text instead of a hard-coded\n// index, so the script keeps working if the table layout shifts. Row 0 is\n// the header row.\nconst headerRow = table.values[0];\nlet eventSizeColumn = -1;\nfor (let colIndex = 0; colIndex < headerRow.length; colIndex++) {\n  if (headerRow[colIndex] === \"Event Size\") {\n    eventSizeColumn = colIndex;\n    break;\n  }\n}\n\nif (eventSizeColumn !== -1) {\n  for (let rowIndex = 1; rowIndex < table.rowCount; rowIndex++) {\n    const cell = table.getCell(rowIndex, eventSizeColumn);\n    const paragraphs = cell.body.paragraphs;\n    paragraphs.load(\"items/text\");\n    await context.sync();\n\n    const paragraph = paragraphs.items[0];\n    if (paragraph.text === \"\") {\n      // Inserting into the paragraph's own (still empty) range reuses its\n      // existing run, so the cell's original formatting is preserved and no\n      // extra paragraph is created.\n      const wholeParagraph = paragraph.getRange(\"Whole\");\n      wholeParagraph.insertText(\"<NA>\", Word.InsertLocation.replace);\n    }\n  }\n  await context.sync();\n}\n", "ps1": "# The document contains one table of structural-variant calls. Two kinds\n# of cleanup are applied to its data cells:\n#\n#   1. Numeric values that were serialized with a trailing \".000\" (e.g.\n#      \"137174039.000\") are rewritten without that suffix (\"137174039\").\n#   2. A few \"Event Size\" cells that were left blank get the literal text\n#      \"<NA>\" inserted.\n\n$d = $word.ActiveDocument\n\n# --- 1. Strip the \".000\" suffix from every matching numeric run -----------\n# Wildcard pattern: one-or-more digits captured in \\1, immediately\n# followed by a literal \".000\" which is dropped on replace.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n[void]$find.Execute(\"([0-9]{1,}).000\", $false, $false, $true, $false, $false, $true, 1, $false, \"\\1\", 2)\n\n# --- 2. Fill blank \"Event Size\" cells with the text \"<NA>\" -----------------\n$t = $d.Tables.Item(1)\n$colCount = $t.Columns.Count\n$rowCount = $t.Rows.Count\n\n# Locate the \"Event Size\" column by its header text instead of a hard-coded\n# index, so the script keeps working if the table layout shifts.\n$eventSizeCol = 0\nfor ($c = 1; $c -le $colCount; $c++) {\n  $headerText = $t.Cell(1, $c).Range.Text -replace \"[\\x07\\x0d]\", \"\"\n  if ($headerText -eq \"Event Size\") {\n    $eventSizeCol = $c\n    break\n  }\n}\n\nif ($eventSizeCol -gt 0) {\n  for ($r = 2; $r -le $rowCount; $r++) {\n    $cell = $t.Cell($r, $eventSizeCol)\n    # A cell's Range.Text always carries trailing cell-mark (0x07) / paragraph\n    # -mark (0x0D) characters; strip those before checking for \"empty\".\n    $clean = $cell.Range.Text -replace \"[\\x07\\x0d]\", \"\"\n    if ($clean -eq \"\") {\n      # Assigning to Range.Text (rather than InsertBefore/TypeText) reuses\n      # the existing run so the cell's original formatting is preserved.\n      $cell.Range.Text = \"<NA>\"\n    }\n  }\n}\n"}
